$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0205992509363296
$ws.Range("C2").Value = 0.917602996254682
$ws.Range("D2").Value = 0.0674157303370786
$ws.Range("E2").Value = 0.925093632958802
$ws.Range("F2").Value = 0.0617977528089888
$ws.Range("G2").Value = 0.0187265917602996
$ws.Range("H2").Value = 0.359550561797753
$ws.Range("I2").Value = 0.0187265917602996
$ws.Range("J2").Value = 0.0205992509363296
$ws.Range("K2").Value = 0.00749063670411985
$ws.Range("L2").Value = 0.0131086142322097
$ws.Range("M2").Value = 0.98876404494382
$ws.Range("N2").Value = 0.00561797752808989
$ws.Range("O2").Value = 0.00936329588014981
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.955056179775281
$ws.Range("R2").Value = 0.00187265917602996
$ws.Range("S2").Value = 0.00561797752808989
$ws.Range("T2").Value = 0.00187265917602996
$ws.Range("U2").Value = 0.00374531835205993
$ws.Range("V2").Value = 0.0842696629213483
$ws.Range("W2").Value = 0.0411985018726592
$ws.Range("X2").Value = 0.99250936329588
$ws.Range("B3").Value = 0.902621722846442
$ws.Range("C3").Value = 0.052434456928839
$ws.Range("D3").Value = 0.0112359550561798
$ws.Range("E3").Value = 0.00749063670411985
$ws.Range("F3").Value = 0.908239700374532
$ws.Range("G3").Value = 0.898876404494382
$ws.Range("H3").Value = 0.629213483146067
$ws.Range("I3").Value = 0.0749063670411985
$ws.Range("J3").Value = 0.0205992509363296
$ws.Range("K3").Value = 0.00936329588014981
$ws.Range("L3").Value = 0.00187265917602996
$ws.Range("M3").Value = 0.00374531835205993
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0.00187265917602996
$ws.Range("P3").Value = 0.932584269662921
$ws.Range("Q3").Value = 0.0411985018726592
$ws.Range("R3").Value = 0.00187265917602996
$ws.Range("S3").Value = 0.0580524344569288
$ws.Range("T3").Value = 0.00187265917602996
$ws.Range("U3").Value = 0.0112359550561798
$ws.Range("V3").Value = 0.00187265917602996
$ws.Range("W3").Value = 0.00187265917602996
$ws.Range("X3").Value = 0.00187265917602996
$ws.Range("B4").Value = 0.00187265917602996
$ws.Range("C4").Value = 0.0168539325842697
$ws.Range("D4").Value = 0.915730337078652
$ws.Range("E4").Value = 0.0617977528089888
$ws.Range("F4").Value = 0.0112359550561798
$ws.Range("G4").Value = 0.00749063670411985
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.00749063670411985
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.179775280898876
$ws.Range("L4").Value = 0.00936329588014981
$ws.Range("M4").Value = 0.00187265917602996
$ws.Range("N4").Value = 0.99250936329588
$ws.Range("O4").Value = 0.00374531835205993
$ws.Range("P4").Value = 0.00187265917602996
$ws.Range("Q4").Value = 0.00374531835205993
$ws.Range("R4").Value = 0.99625468164794
$ws.Range("S4").Value = 0.00936329588014981
$ws.Range("T4").Value = 0.99625468164794
$ws.Range("U4").Value = 0.00749063670411985
$ws.Range("V4").Value = 0.908239700374532
$ws.Range("W4").Value = 0.955056179775281
$ws.Range("X4").Value = 0.00374531835205993
$ws.Range("B5").Value = 0.0730337078651685
$ws.Range("C5").Value = 0.0131086142322097
$ws.Range("D5").Value = 0.00561797752808989
$ws.Range("E5").Value = 0.00561797752808989
$ws.Range("F5").Value = 0.0187265917602996
$ws.Range("G5").Value = 0.0749063670411985
$ws.Range("H5").Value = 0.0112359550561798
$ws.Range("I5").Value = 0.898876404494382
$ws.Range("J5").Value = 0.958801498127341
$ws.Range("K5").Value = 0.803370786516854
$ws.Range("L5").Value = 0.975655430711611
$ws.Range("M5").Value = 0.00561797752808989
$ws.Range("N5").Value = 0.00187265917602996
$ws.Range("O5").Value = 0.98501872659176
$ws.Range("P5").Value = 0.0655430711610487
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.926966292134832
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0.97752808988764
$ws.Range("V5").Value = 0.00374531835205993
$ws.Range("W5").Value = 0.00187265917602996
$ws.Range("X5").Value = 0.00187265917602996

Write-Host "Updated frequency table values for gRNA-9 run"
